$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Update the confidential disclaimer text in A41 (date 2021-04-22 -> 2021-04-23)
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-38
$ws.Range("D2").Value = 0.030288518435951
$ws.Range("E2").Value = 0.01208541572012711
$ws.Range("D3").Value = 0.02991834072312634
$ws.Range("E3").Value = 0.01685426481893981
$ws.Range("D4").Value = 0.02953046997038972
$ws.Range("E4").Value = 0.04679989881102964
$ws.Range("D5").Value = 0.06505219643363361
$ws.Range("E5").Value = 0.009622126054686708
$ws.Range("D6").Value = 0.01655242860563512
$ws.Range("E6").Value = -0.05857621321171513
$ws.Range("D7").Value = 0.01568429011395497
$ws.Range("E7").Value = 0.01850041362713406
$ws.Range("D8").Value = 0.0300840655303017
$ws.Range("E8").Value = 0.007763183689472486
$ws.Range("D9").Value = 0.03261810202435899
$ws.Range("E9").Value = 0.05135004821600786
$ws.Range("D10").Value = 0.03000621615468907
$ws.Range("E10").Value = -0.006348519988993218
$ws.Range("D11").Value = 0.0315614343629503
$ws.Range("E11").Value = 0.01190320470896022
$ws.Range("D12").Value = 0.01318524652105688
$ws.Range("E12").Value = 0.03578350976591649
$ws.Range("D13").Value = 0.01482086976625135
$ws.Range("E13").Value = 0.01132776230269283
$ws.Range("D14").Value = 0.01620938022067546
$ws.Range("E14").Value = 0.0191866881738687
$ws.Range("D15").Value = 0.007824255427733172
$ws.Range("E15").Value = 0.03356783919597994
$ws.Range("D16").Value = 0.007055984316889422
$ws.Range("E16").Value = 0.02006018054162495
$ws.Range("D17").Value = 0.03183665942824745
$ws.Range("E17").Value = 0.009336503133779983
$ws.Range("D18").Value = 0.03033963166236333
$ws.Range("E18").Value = 0.001282965074839604
$ws.Range("D19").Value = 0.03136386208393342
$ws.Range("E19").Value = 0.006111320045129709
$ws.Range("D20").Value = 0.02914633441496784
$ws.Range("E20").Value = 0.01554701200593556
$ws.Range("D21").Value = 0.04458488786347893
$ws.Range("E21").Value = 0.02089138946700064
$ws.Range("D22").Value = 0.03252079030484321
$ws.Range("E22").Value = 0.02569140093698041
$ws.Range("D23").Value = 0.031597803389436
$ws.Range("E23").Value = 0.007559260872270102
$ws.Range("D24").Value = 0.030005626386692
$ws.Range("E24").Value = -0.002299663895276804
$ws.Range("D25").Value = 0.01489400099788745
$ws.Range("E25").Value = 0.01936326918507936
$ws.Range("D26").Value = 0.01467106869499676
$ws.Range("E26").Value = -0.001849171892587043
$ws.Range("D27").Value = 0.03055155496264213
$ws.Range("E27").Value = 0.008139864099660388
$ws.Range("D28").Value = 0.0307198354311381
$ws.Range("E28").Value = 0.003942046792607234
$ws.Range("D29").Value = 0.03033412716105739
$ws.Range("E29").Value = 0.01547614418478038
$ws.Range("D30").Value = 0.02793495094899568
$ws.Range("E30").Value = 0.007818547762811434
$ws.Range("D31").Value = 0.03503280879367674
$ws.Range("E31").Value = 0.0279456574805137
$ws.Range("D32").Value = 0.03158738415482118
$ws.Range("E32").Value = -0.001643047853768675
$ws.Range("D33").Value = 0.02947188634934789
$ws.Range("E33").Value = 0.02878945542837319
$ws.Range("D34").Value = 0.03118221354083731
$ws.Range("E34").Value = 0.009507237605285601
$ws.Range("D35").Value = 0.0305311096720772
$ws.Range("E35").Value = -0.0002318034306907357
$ws.Range("D36").Value = 0.0289713699091718
$ws.Range("E36").Value = 0.01913550926240082
$ws.Range("D37").Value = 0.03233029524179112
$ws.Range("E37").Value = 0.002140390134747205
$ws.Range("E38").Value = 0.01257267416143848

$ws.Protect("D382")
